# Updates the cryptos price/volume table on Sheet1 (rows 2-51) to the
# latest scraped values. Columns B (Coin) and C (Link) are only touched
# for the two rows whose ranking swapped (43/44: WhiteBITCoin <-> dogwifhat).
# Column D (Price) holds text, not numbers (e.g. "68.232.90" uses dots as
# thousands separators and isn't a valid float) -- for values that DO look
# like a plain decimal (e.g. "599.47"), Excel's Range.Value setter would
# silently coerce the cell to a Number, so we force the cell to Text via
# NumberFormat "@" before assigning, then restore General/Normal so the
# cell keeps its original (unstyled) look while the stored value stays a
# string, matching the workbook's inlineStr-everywhere convention.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '68.287.82'
$ws.Range("E2").Value = '  +1.38%  '

# Row 3
$ws.Range("D3").Value = '2.639.28'
$ws.Range("E3").Value = '  +1.22%  '

# Row 4
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '599.47'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.37%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.60'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.76%  '

# Row 7
$ws.Range("E7").Value = '  -0.03%  '

# Row 8
$ws.Range("E8").Value = '  +0.01%  '

# Row 9
$ws.Range("D9").Value = '2.637.70'
$ws.Range("E9").Value = '  +1.26%  '

# Row 10
$ws.Range("E10").Value = '  +7.30%  '

# Row 11
$ws.Range("E11").Value = '  -0.67%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.24'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.26%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.350'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.74%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.02'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.73%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000188'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.28%  '

# Row 16
$ws.Range("D16").Value = '3.114.15'
$ws.Range("E16").Value = '  +1.06%  '

# Row 17
$ws.Range("D17").Value = '68.195.93'
$ws.Range("E17").Value = '  +1.49%  '

# Row 18
$ws.Range("D18").Value = '2.634.52'
$ws.Range("E18").Value = '  +0.90%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.43'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.56%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '366.19'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.46%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.41'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.71%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.28'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.03%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.85'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.31%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.12'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.18%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '73.63'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.45%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.06%  '

# Row 27
$ws.Range("E27").Value = '  +1.31%  '

# Row 28
$ws.Range("D28").Value = '2.768.16'
$ws.Range("E28").Value = '  +1.20%  '

# Row 29
$ws.Range("E29").Value = '  +5.96%  '

# Row 30
$ws.Range("E30").Value = '  -0.15%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '575.25'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.40%  '

# Row 32
$ws.Range("E32").Value = '  +4.33%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.01'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.37%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.86'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.53%  '

# Row 35
$ws.Range("E35").Value = '  +3.26%  '

# Row 36
$ws.Range("E36").Value = '  -0.03%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.55'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.56%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '160.16'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.33%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.28'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.13%  '

# Row 40
$ws.Range("E40").Value = '  +3.47%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.369'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.93%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.39'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.05%  '

# Row 43
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.64'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.70%  '

# Row 44
$ws.Range("B44").Value = 'WhiteBITCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.73'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.50%  '

# Row 45
$ws.Range("D45").Value = '0.0₆0321'
$ws.Range("E45").Value = '  +13.61%  '

# Row 46
$ws.Range("E46").Value = '  +0.03%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '40.51'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.32%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '157.23'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.58%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.77'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.61%  '

# Row 50
$ws.Range("E50").Value = '  +1.91%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '21.91'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.71%  '

